$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (the "Förändrad" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45210) {
        $cell.Value2 = 45211
    }
}
